# SwaadSutra_Daily_2026-01-20.xlsx update
# New order (#23) placed at 2026-01-20 12:17 by Radha shelke -> inserted as the
# newest row at the top of the "Daily Orders" log, pushing existing orders down.
# Summary totals and the Items Breakdown pivot are updated to reflect the new order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Daily Orders - insert the new order as row 2 (existing rows shift down)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Orders")

$ws1.Rows("2:2").Insert()

$ws1.Range("A2").Value = 23
$ws1.Range("B2").Value = "2026-01-20 12:17"
$ws1.Range("C2").Value = "Radha shelke"
$ws1.Range("D2").Value = "C 803"
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "9890774770"
$ws1.Range("F2").Value = "Appe Chutney x2, Onion Pakoda (Kanda Bhaje) x1"
$ws1.Range("G2").Value = 180
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"
# Collection Date / Collection Time / Notes / Cancel Reason / Feedback stay blank
# for a brand new order (matches the existing blank cells on the other rows).

# ---------------------------------------------------------------------------
# Sheet 2: Summary - bump Total Orders / New / Total Revenue
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("A2").Value = 4
$ws2.Range("B2").Value = 1
$ws2.Range("G2").Value = 260

# ---------------------------------------------------------------------------
# Sheet 3: Items Breakdown - add the two new items, keeping the existing
# quantity-desc / revenue-desc ordering used by the report.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Items Breakdown")

# "Appe Chutney" (qty 2, revenue 120) becomes the new top row.
$ws3.Rows("2:2").Insert()
$ws3.Range("A2").Value = "Appe Chutney"
$ws3.Range("B2").Value = 2
$ws3.Range("C2").Value = 120

# "Onion Pakoda (Kanda Bhaje)" (qty 1, revenue 60) slots in ahead of the other
# qty-1 items (after Wheat Chapati, which is now row 3).
$ws3.Rows("4:4").Insert()
$ws3.Range("A4").Value = "Onion Pakoda (Kanda Bhaje)"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 60
